$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.919.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.696.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +9.30%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.684.24"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.09%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.203"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.612"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000287"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.289.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "688.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("E16").Value = "  +4.50%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.713.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +10.02%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "72.047.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.68%  "
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.944"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.81%  "
$ws.Range("E24").Value = "  +3.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.95%  "
$ws.Range("E29").Value = "  +6.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "579.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("E35").Value = "  +3.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "60.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.752.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("E40").Value = "  +6.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.77%  "
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0461"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.348"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("E47").Value = "  +6.95%  "
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("E49").Value = "  +4.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.66%  "
